$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "66.605.54"
$ws.Cells.Item(2, 5).Value = "  +0.40%  "
$ws.Cells.Item(3, 4).Value = "3.561.47"
$ws.Cells.Item(3, 5).Value = "  +0.67%  "
$ws.Cells.Item(4, 5).Value = "  -0.03%  "
$ws.Cells.Item(5, 4).Value = "'607.43"
$ws.Cells.Item(5, 5).Value = "  -0.18%  "
$ws.Cells.Item(6, 4).Value = "'145.17"
$ws.Cells.Item(6, 5).Value = "  +0.64%  "
$ws.Cells.Item(7, 4).Value = "3.559.85"
$ws.Cells.Item(7, 5).Value = "  +0.67%  "
$ws.Cells.Item(8, 5).Value = "  -0.05%  "
$ws.Cells.Item(9, 4).Value = "'0.497"
$ws.Cells.Item(9, 5).Value = "  +3.51%  "
$ws.Cells.Item(10, 5).Value = "  -0.64%  "
$ws.Cells.Item(11, 4).Value = "'7.98"
$ws.Cells.Item(11, 5).Value = "  -2.09%  "
$ws.Cells.Item(12, 4).Value = "'0.415"
$ws.Cells.Item(12, 5).Value = "  +0.87%  "
$ws.Cells.Item(13, 4).Value = "4.164.00"
$ws.Cells.Item(13, 5).Value = "  +0.62%  "
$ws.Cells.Item(14, 5).Value = "  -0.14%  "
$ws.Cells.Item(15, 4).Value = "'30.06"
$ws.Cells.Item(15, 5).Value = "  -0.74%  "
$ws.Cells.Item(16, 4).Value = "3.557.73"
$ws.Cells.Item(16, 5).Value = "  +0.59%  "
$ws.Cells.Item(17, 4).Value = "66.624.58"
$ws.Cells.Item(17, 5).Value = "  +0.32%  "
$ws.Cells.Item(18, 5).Value = "  +0.34%  "
$ws.Cells.Item(19, 5).Value = "  +5.08%  "
$ws.Cells.Item(20, 4).Value = "'6.22"
$ws.Cells.Item(20, 5).Value = "  +0.10%  "
$ws.Cells.Item(21, 4).Value = "'14.95"
$ws.Cells.Item(21, 5).Value = "  +0.03%  "
$ws.Cells.Item(22, 4).Value = "'431.77"
$ws.Cells.Item(22, 5).Value = "  +1.30%  "
$ws.Cells.Item(23, 5).Value = "  +2.27%  "
$ws.Cells.Item(24, 4).Value = "'79.48"
$ws.Cells.Item(24, 5).Value = "  +0.80%  "
$ws.Cells.Item(25, 4).Value = "3.701.39"
$ws.Cells.Item(25, 5).Value = "  +0.71%  "
$ws.Cells.Item(26, 4).Value = "'1.00"
$ws.Cells.Item(26, 5).Value = "  -0.18%  "
$ws.Cells.Item(27, 5).Value = "  -0.50%  "
$ws.Cells.Item(28, 4).Value = "'8.01"
$ws.Cells.Item(28, 5).Value = "  -1.99%  "
$ws.Cells.Item(29, 5).Value = "  +1.14%  "
$ws.Cells.Item(30, 5).Value = "  -1.16%  "
$ws.Cells.Item(31, 5).Value = "  -0.02%  "
$ws.Cells.Item(32, 5).Value = "  -2.36%  "
$ws.Cells.Item(33, 4).Value = "3.555.23"
$ws.Cells.Item(33, 5).Value = "  +0.76%  "
$ws.Cells.Item(34, 5).Value = "  +0.25%  "
$ws.Cells.Item(35, 5).Value = "  -3.96%  "
$ws.Cells.Item(36, 5).Value = "  +0.44%  "
$ws.Cells.Item(37, 5).Value = "  +0.02%  "
$ws.Cells.Item(38, 5).Value = "  -1.77%  "
$ws.Cells.Item(39, 5).Value = "  -0.28%  "
$ws.Cells.Item(40, 4).Value = "'173.92"
$ws.Cells.Item(40, 5).Value = "  +0.21%  "
$ws.Cells.Item(41, 4).Value = "'0.0849"
$ws.Cells.Item(41, 5).Value = "  -0.97%  "
$ws.Cells.Item(42, 5).Value = "  -1.49%  "
$ws.Cells.Item(43, 5).Value = "  -0.70%  "
$ws.Cells.Item(44, 5).Value = "  +2.06%  "
$ws.Cells.Item(45, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(45, 4).Value = "'1.00"
$ws.Cells.Item(45, 5).Value = "  +0.04%  "
$ws.Cells.Item(46, 2).Value = "dogwifhat"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(46, 4).Value = "'2.53"
$ws.Cells.Item(46, 5).Value = "  +5.19%  "
$ws.Cells.Item(47, 2).Value = "ONDO"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Cells.Item(47, 4).Value = "'1.18"
$ws.Cells.Item(47, 5).Value = "  -3.01%  "
$ws.Cells.Item(48, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(48, 4).Value = "'25.10"
$ws.Cells.Item(48, 5).Value = "  -3.91%  "
$ws.Cells.Item(49, 2).Value = "Cosmos"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(49, 4).Value = "'7.17"
$ws.Cells.Item(49, 5).Value = "  +0.46%  "
$ws.Cells.Item(50, 2).Value = "EnergySwap"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(50, 4).Value = "'23.54"
$ws.Cells.Item(50, 5).Value = "  +4.42%  "
$ws.Cells.Item(51, 2).Value = "SuiNetwork"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Cells.Item(51, 4).Value = "'0.942"
$ws.Cells.Item(51, 5).Value = "  -0.29%  "
